$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Sheet1" to "Sheet"
$ws.Name = "Sheet"

# Clear column C (formulas no longer needed)
$ws.Range("C1:C2").Clear()

# Replace numeric data with text labels / values
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Value"
$ws.Range("A2").Value = "Test"

# B2 holds the digit string "123" as TEXT, not a number -- force a text
# number format first so the value is stored as a string, not coerced.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "123"
